$d = $word.ActiveDocument

# Locate the title text that needs to change ("CS202 Chapter 3" -> "CS202 Chapter 4").
$titleRange = $d.Content.Duplicate
$titleRange.Find.Execute("CS202 Chapter 3", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Word moves its "_GoBack" bookmark to the location of the last edit. Re-adding a
# bookmark with the same name automatically removes the previous one, so placing it
# right after the edited title text reproduces that move (it currently sits around
# the two inline images further down the document).
$editEnd = $d.Range($titleRange.End, $titleRange.End)
$d.Bookmarks.Add("_GoBack", $editEnd)

# Update the title text itself, keeping the edit scoped to the matched range so the
# following " Homework" run is left untouched.
$titleRange.Text = "CS202 Chapter 4"
